# Effort Calculation workbook - adding effort entries for 26.10.2015
# (Business: "Technical Template creation"; School: "Discussion on how to
# progress with requirement") plus normalising the School sheet's existing
# date cells to the same text format used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook

# --- Business sheet: new entry on row 4 -----------------------------------
$business = $wb.Worksheets.Item("Business")
$business.Range("A4").Value = "26.10.2015"
$business.Range("B4").Value = "Arpan Kar"
$business.Range("C4").Value = "Architecture"
$business.Range("D4").Value = 1
$business.Range("J4").Value = "Technical Template creation"

# --- School sheet: normalise earlier date cells + new entry on row 5 ------
$school = $wb.Worksheets.Item("School")
$school.Range("A2").Value = "21.10.2015"
$school.Range("A3").Value = "22.10.2015"
$school.Range("A4").Value = "24.10.2015"

$school.Range("A5").Value = "26.10.2015"
$school.Range("B5").Value = "Arpan Kar"
$school.Range("C5").Value = "Architecture"
$school.Range("D5").Value = 0.5
$school.Range("J5").Value = "Discussion on how to progress with requirement"
$school.Range("A5").Select()

# --- Team Member sheet: trim the lingering wide selection ------------------
$teamMember = $wb.Worksheets.Item("Team Member")
$teamMember.Select()
$teamMember.Range("C5").Select()

# --- Service sheet was the previously active tab; Cover becomes active now -
$service = $wb.Worksheets.Item("Service")
$service.Select()

# --- Cover is the final active sheet ---------------------------------------
$cover = $wb.Worksheets.Item("Cover")
$cover.Activate()
$cover.Range("A1").Select()
